$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added 64 bit designs: refreshed power report numbers in row 2.
$ws.Range("B2").Value = 0.015711084008216858
$ws.Range("C2").Value = 0.005493159871548414
$ws.Range("D2").Value = 0.0040701813995838165
$ws.Range("E2").Value = 0.003394484054297209
$ws.Range("F2").Value = 0.0000000392246626290671
$ws.Range("G2").Value = 0.001065471675246954
$ws.Range("J2").Value = 0.12726984918117523
$ws.Range("K2").Value = 1.415169358253479

# Column G (XADC) narrowed to match the other 8.75-wide columns.
$ws.Columns.Item(7).ColumnWidth = 8
